$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "C7" = 12014
    "D7" = 12016
    "E7" = 12013
    "F7" = 12015
    "G7" = 12019
    "H7" = 12017
    "C8" = 12020
    "D8" = 12018
    "E8" = 12022
    "F8" = 12024
    "G8" = 12021
    "H8" = 12023
    "C12" = 12016
    "D12" = 12013
    "E12" = 12026
    "F12" = 12028
    "G12" = 12025
    "H12" = 12027
    "C13" = 12021
    "D13" = 12019
    "E13" = 12001
    "F13" = 12023
    "G13" = 12020
    "H13" = 12022
    "D14" = 12024
    "E14" = 12004
    "F14" = 12006
    "G14" = 12003
    "H14" = 12005
    "C15" = 12008
    "D15" = 12010
    "E15" = 12007
    "F15" = 12011
    "G15" = 12012
    "H15" = 12009
    "C16" = 12026
    "D16" = 12027
    "E16" = 12019
    "F16" = 12001
    "G16" = 12028
    "H16" = 12025
    "C19" = 12004
    "D19" = 12005
    "E19" = 12002
    "F19" = 12003
    "G19" = 12006
    "H19" = 12007
    "C20" = 12012
    "D20" = 12001
    "E20" = 12008
    "F20" = 12009
    "G20" = 12010
    "H20" = 12011
    "C21" = 12013
    "D21" = 12017
    "E21" = 12016
    "F21" = 12018
    "G21" = 12014
    "H21" = 12015
    "C22" = 12027
    "D22" = 12004
    "E22" = 12002
    "F22" = 12025
    "H22" = 12003
    "C23" = 12001
    "D23" = 12006
    "E23" = 12028
    "F23" = 12005
    "G23" = 12007
    "H23" = 12010
    "C24" = 12008
    "D24" = 12011
    "E24" = 12012
    "F24" = 12014
    "G24" = 12009
    "H24" = 12013
    "C27" = 12022
    "D27" = 12023
    "E27" = 12020
    "F27" = 12021
    "G27" = 12024
    "H27" = 12019
    "D28" = 12015
    "E28" = 12018
    "F28" = 12003
    "G28" = 12002
    "H28" = 12016
    "C29" = 12006
    "D29" = 12011
    "E29" = 12005
    "F29" = 12008
    "H29" = 12010
    "C30" = 12015
    "D30" = 12022
    "E30" = 12012
    "F30" = 12020
    "G30" = 12014
    "H30" = 12021
    "E31" = 12023
    "F31" = 12017
    "G31" = 12018
    "H31" = 12024
    "C32" = 12025
    "D32" = 12026
    "E32" = 12027
    "F32" = 12028
}

foreach ($addr in $changes.Keys) {
    $ws.Range($addr).Value = $changes[$addr]
}

# The last row (row 33, match 31) is removed entirely - the judging
# session time between rounds changed, eliminating the final round.
$ws.Rows.Item(33).Delete()
